# Day 16 PPT.pptx edit — resize title placeholders and reposition a picture.
#
# Slide order in this deck (Slides.Item index) vs. XML part:
#   3  -> slide3.xml  (title "Data Structures Introduction")
#   9  -> slide9.xml  (title "Java Collection Framework")
#  11  -> slide11.xml (picture "Picture 5")
#  13  -> slide13.xml (title "Iterating Over a Collection")
#  17  -> slide17.xml (title "Maps")
#  22  -> slide22.xml (title "Comparing Custom Objects")

$p = $ppt.ActivePresentation

# --- Slide 3: "Data Structures Introduction" title placeholder — shrink height ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$sh3.Height = 269.9349060058594

# --- Slide 9: "Java Collection Framework" title placeholder — shrink height ---
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(2)
$sh9.Height = 268.891357421875

# --- Slide 11: "Picture 5" — reposition slightly ---
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(4)
$sh11.Left = 13.718740157480315
$sh11.Top = 110.90772247314453

# --- Slide 13: "Iterating Over a Collection" title placeholder — shrink height ---
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$sh13.Height = 269.9349060058594

# --- Slide 17: "Maps" title placeholder — shrink height ---
$s17 = $p.Slides.Item(17)
$sh17 = $s17.Shapes.Item(2)
$sh17.Height = 277.2392125984252

# --- Slide 22: "Comparing Custom Objects" title placeholder — shrink height ---
$s22 = $p.Slides.Item(22)
$sh22 = $s22.Shapes.Item(2)
$sh22.Height = 274.1087646484375
